# Update cryptocurrency price and volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.051.77"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").Value = "2.987.24"
$ws.Range("E3").Value = "  -2.09%  "

$ws.Range("D5").Value = "'501.77"
$ws.Range("E5").Value = "  -4.51%  "

$ws.Range("D6").Value = "'138.19"
$ws.Range("E6").Value = "  -3.08%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -3.17%  "

$ws.Range("D9").Value = "'7.31"
$ws.Range("E9").Value = "  -4.16%  "

$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("D11").Value = "'0.359"
$ws.Range("E11").Value = "  -2.26%  "

$ws.Range("D12").Value = "3.492.28"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("D14").Value = "'26.07"
$ws.Range("E14").Value = "  -1.52%  "

$ws.Range("E15").Value = "  -4.59%  "

$ws.Range("D16").Value = "57.121.82"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "2.992.78"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").Value = "'12.64"
$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("D20").Value = "'7.88"
$ws.Range("E20").Value = "  -3.67%  "

$ws.Range("D21").Value = "'321.38"
$ws.Range("E21").Value = "  -5.71%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("D25").Value = "'63.82"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").Value = "'0.164"
$ws.Range("E27").Value = "  -5.36%  "

$ws.Range("D28").Value = "0.0₃0897"
$ws.Range("E28").Value = "  -7.02%  "

$ws.Range("E29").Value = "  -4.73%  "

$ws.Range("D30").Value = "'7.07"
$ws.Range("E30").Value = "  -2.68%  "

$ws.Range("D31").Value = "'1.78"
$ws.Range("E31").Value = "  -4.40%  "

$ws.Range("E32").Value = "  -5.40%  "

$ws.Range("D33").Value = "'20.22"
$ws.Range("E33").Value = "  -4.00%  "

$ws.Range("D34").Value = "'155.06"
$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("D35").Value = "'4.58"
$ws.Range("E35").Value = "  -2.97%  "

$ws.Range("D36").Value = "'5.78"
$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("E37").Value = "  -6.15%  "

$ws.Range("D38").Value = "'24.25"
$ws.Range("E38").Value = "  -6.63%  "

$ws.Range("D39").Value = "'0.0666"
$ws.Range("E39").Value = "  -3.61%  "

$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("D41").Value = "3.017.54"
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "'3.75"
$ws.Range("E43").Value = "  -2.50%  "

$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("D45").Value = "2.198.07"

$ws.Range("E46").Value = "  -5.92%  "

$ws.Range("E47").Value = "  -0.94%  "

$ws.Range("D48").Value = "'0.939"
$ws.Range("E48").Value = "  -8.81%  "

$ws.Range("E49").Value = "  -4.64%  "

$ws.Range("D50").Value = "'19.34"
$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").Value = "'1.82"
$ws.Range("E51").Value = "  -10.79%  "
